# Update stock report figures (quantities, values and sub/grand totals)
# to reflect corrected company-wise stock data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F41").Value = 200
$ws.Range("G41").Value = 38578
$ws.Range("B66").Value = 186930.34
$ws.Range("F149").Value = 50
$ws.Range("G149").Value = 9064.5
$ws.Range("B155").Value = 34785.01
$ws.Range("F215").Value = 161
$ws.Range("G215").Value = 18078.69
$ws.Range("B218").Value = 70279.33
$ws.Range("F222").Value = 570
$ws.Range("G222").Value = 10545
$ws.Range("B229").Value = 19497.2
$ws.Range("F266").Value = 7
$ws.Range("G266").Value = 466.13
$ws.Range("F286").Value = 18
$ws.Range("G286").Value = 1573.02
$ws.Range("B295").Value = 106580.07
$ws.Range("B304").Value = 63520
$ws.Range("E304").Value = 153.4
$ws.Range("F304").Value = 36
$ws.Range("G304").Value = 5194.08
$ws.Range("B305").Value = 55373
$ws.Range("E305").Value = 163.62
$ws.Range("F305").Value = -94
$ws.Range("G305").Value = -13562.32
$ws.Range("B306").Value = 63531
$ws.Range("E306").Value = 152.53
$ws.Range("F306").Value = 26
$ws.Range("G306").Value = 3730.48
$ws.Range("B307").Value = 57802
$ws.Range("E307").Value = 162.71
$ws.Range("F307").Value = -79
$ws.Range("G307").Value = -11334.92
$ws.Range("B308").Value = 63510
$ws.Range("E308").Value = 50.66
$ws.Range("F308").Value = 76
$ws.Range("G308").Value = 3620.64
$ws.Range("B309").Value = 55356
$ws.Range("E309").Value = 54.04
$ws.Range("F309").Value = -158
$ws.Range("G309").Value = -7527.12
$ws.Range("B317").Value = 60325
$ws.Range("E317").Value = 151.57
$ws.Range("F317").Value = -102
$ws.Range("G317").Value = -12939.72
$ws.Range("B318").Value = 63560
$ws.Range("E318").Value = 134.87
$ws.Range("F318").Value = 1
$ws.Range("G318").Value = 126.86
$ws.Range("B381").Value = 58047
$ws.Range("D381").Value = 105.54
$ws.Range("E381").Value = 126.1
$ws.Range("F381").Value = 32
$ws.Range("G381").Value = 3377.28
$ws.Range("B382").Value = 47097
$ws.Range("D382").Value = 112.28
$ws.Range("E382").Value = 134.16
$ws.Range("F382").Value = 15
$ws.Range("G382").Value = 1684.2
$ws.Range("F430").Value = 215
$ws.Range("G430").Value = 9950.200000000001
$ws.Range("B438").Value = 22723.02
$ws.Range("B479").Value = 53319
$ws.Range("E479").Value = 310.64
$ws.Range("F479").Value = -6
$ws.Range("G479").Value = -1643.52
$ws.Range("B480").Value = 64810
$ws.Range("E480").Value = 291.22
$ws.Range("F480").Value = 0
$ws.Range("G480").Value = 0
$ws.Range("F488").Value = 2
$ws.Range("G488").Value = 127.22
$ws.Range("B493").Value = 3121.97
$ws.Range("B496").Value = 60025
$ws.Range("E496").Value = 37.22
$ws.Range("F496").Value = -98
$ws.Range("G496").Value = -3217.34
$ws.Range("B497").Value = 64833
$ws.Range("E497").Value = 34.9
$ws.Range("F497").Value = 88
$ws.Range("G497").Value = 2889.04
$ws.Range("F498").Value = 127
$ws.Range("G498").Value = 4169.41
$ws.Range("B506").Value = 64830
$ws.Range("E506").Value = 34.9
$ws.Range("F506").Value = 83
$ws.Range("G506").Value = 2724.89
$ws.Range("B507").Value = 60022
$ws.Range("E507").Value = 37.22
$ws.Range("F507").Value = -113
$ws.Range("G507").Value = -3709.79
$ws.Range("B508").Value = 41090.9
$ws.Range("F610").Value = 58
$ws.Range("G610").Value = 1469.14
$ws.Range("B628").Value = 199951.49
$ws.Range("F646").Value = 3
$ws.Range("G646").Value = 2849.82
$ws.Range("F647").Value = 1
$ws.Range("G647").Value = 949.9400000000001
$ws.Range("B657").Value = 68002.7
$ws.Range("B718").Value = 2425977.62
$ws.Range("B719").Value = 2425977.62
